$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("n" = number of processors): drop from 1,000,000 down to 50 for
# every data column so the non-optimized (2^(n/2)) estimate stays finite.
$ws.Range("B6:I6").Value = 50

# Row 7 ("lambda"): the per-processor execution-time estimate now divides
# by the extra 2^(B6/2) factor that models the un-optimized algorithm's
# blow-up. Only the master cell needs updating -- C7 (=B7) and the shared
# D7:I7 group recompute automatically.
$ws.Range("B7").Formula = "=B2/(B6 * 2^(B6 / 2))"
# Setting .Formula on a General-formatted cell that references B2 (which
# carries the 0.00000 custom number format) makes Excel copy that number
# format onto B7. The source workbook keeps B7 on the default "Normal"
# style, so restore it explicitly.
$ws.Range("B7").Style = "Normal"

# Row 3 ("Estimate execution time"): rebuild the formula around the new
# lambda definition -- B7 * B6 * 2^(B6/2) / B1 -- instead of the old
# LOG-based one. B3 is its own formula; C3:I3 share one formula (si="0").
$ws.Range("B3").Formula = "=B7 * B6 * 2^ (B6 / 2) / B1"
$ws.Range("C3:I3").Formula = "=C7 * C6 * 2^ (C6 / 2) / C1"

# Reflect the author's final selection (row 3, columns B:I) when the file
# is reopened.
$ws.Range("B3:I3").Select()
